$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "98.352.97"

$ws.Range("D3").Value = "3.420.68"
$ws.Range("E3").Value = "  +0.71%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").Value = "'255.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.28%  "

$ws.Range("D6").Value = "'686.26"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.93%  "

$ws.Range("D7").Value = "'1.46"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -6.46%  "

$ws.Range("D8").Value = "'0.435"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -4.42%  "

$ws.Range("D9").Value = "'1.06"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.45%  "

$ws.Range("E10").Value = "  +0.05%  "

$ws.Range("D11").Value = "3.416.08"
$ws.Range("E11").Value = "  +0.65%  "

$ws.Range("E12").Value = "  +3.56%  "

$ws.Range("D13").Value = "'42.03"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.49%  "

$ws.Range("D14").Value = "'6.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +13.41%  "

$ws.Range("D15").Value = "98.031.85"
$ws.Range("E15").Value = "  -0.42%  "

$ws.Range("E16").Value = "  +0.57%  "

$ws.Range("D17").Value = "4.053.25"
$ws.Range("E17").Value = "  +1.19%  "

$ws.Range("D18").Value = "'9.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +17.67%  "

$ws.Range("D19").Value = "3.423.09"
$ws.Range("E19").Value = "  +1.73%  "

$ws.Range("D20").Value = "'0.582"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +33.53%  "

$ws.Range("D21").Value = "'17.69"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.21%  "

$ws.Range("D22").Value = "'11.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.39%  "

$ws.Range("D23").Value = "'3.45"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.23%  "

$ws.Range("D24").Value = "'511.34"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.33%  "

$ws.Range("D25").Value = "'0.0000207"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.72%  "

$ws.Range("D26").Value = "'6.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.13%  "

$ws.Range("D27").Value = "'101.50"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.21%  "

$ws.Range("D28").Value = "'12.86"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.90%  "

$ws.Range("D29").Value = "3.607.81"
$ws.Range("E29").Value = "  +1.88%  "

$ws.Range("D31").Value = "'11.63"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.93%  "

$ws.Range("D32").Value = "'0.198"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.33%  "

$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.04%  "

$ws.Range("D34").Value = "'2.70"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +25.47%  "

$ws.Range("D35").Value = "'0.578"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.16%  "

$ws.Range("E36").Value = "  +0.25%  "

$ws.Range("D37").Value = "'30.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.30%  "

$ws.Range("D38").Value = "'8.09"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.63%  "

$ws.Range("D39").Value = "'1.54"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +13.77%  "

$ws.Range("D40").Value = "'537.75"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.26%  "

$ws.Range("D41").Value = "'0.154"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.47%  "

$ws.Range("E42").Value = "  -0.02%  "

$ws.Range("D43").Value = "'0.880"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.03%  "

$ws.Range("D44").Value = "'24.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").Value = "'0.0439"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -3.70%  "

$ws.Range("B46").Value = "Cosmos"
$ws.Range("C46").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D46").Value = "'9.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +14.54%  "

$ws.Range("B47").Value = "ImmutableX"
$ws.Range("C47").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D47").Value = "'1.77"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +15.19%  "

$ws.Range("B48").Value = "MantraDAO"
$ws.Range("C48").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D48").Value = "'3.80"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.17%  "

$ws.Range("B49").Value = "Filecoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D49").Value = "'5.84"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +14.43%  "

$ws.Range("B50").Value = "OKB"
$ws.Range("C50").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D50").Value = "'55.84"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.64%  "

$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "'3.27"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.19%  "

